$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166417956352234
$ws.Range("B1").Value = 2.429310321807861
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.370743751525879
$ws.Range("E1").Value = 1.234593987464905
